$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Home Summary
# ---------------------------------------------------------------------------
$home = $wb.Worksheets.Item("Home Summary")

$home.Range("B5").Value  = "KES 886,450"
$home.Range("B6").Value  = "KES 413,550"
$home.Range("B7").Value  = "68.19%"
$home.Range("B8").Value  = "KES 9,963"

$home.Range("B12").Value = "KES 46,100"
$home.Range("B13").Value = "KES 76,100"
$home.Range("B14").Value = "KES 962,550"
$home.Range("B15").Value = "74.04%"
$home.Range("B16").Value = "KES 337,450"

$home.Range("B25").Value = 24160
$home.Range("C25").Value = 250
$home.Range("D25").Value = 24410
$home.Range("E25").Value = "1.88%"

# ---------------------------------------------------------------------------
# Sheet: Daily Expenses
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily Expenses")

# Shift old rows 297-301 (the 03/10/2025 labor entries) down to 299-303.
# Copy bottom-up so the source data isn't clobbered before it is read.
$daily.Range("A301:I301").Copy($daily.Range("A303:I303"))
$daily.Range("A300:I300").Copy($daily.Range("A302:I302"))
$daily.Range("A299:I299").Copy($daily.Range("A301:I301"))
$daily.Range("A298:I298").Copy($daily.Range("A300:I300"))
$daily.Range("A297:I297").Copy($daily.Range("A299:I299"))

# Rows 297 & 298 become new "Workers Accommodation / Food" entries.
# Seed them from the plainly-styled row 2 template so they pick up style s4
# (the same look already used by ordinary, non-pending rows) without minting
# any brand new style entries.
$daily.Range("A2:I2").Copy($daily.Range("A297:I297"))
$daily.Range("A2:I2").Copy($daily.Range("A298:I298"))

# Fix up the date column by copying the already-existing "02/10/2025" text
# value (values-only, so the destination keeps its style and the date stays
# plain text instead of being reinterpreted as a date serial).
$daily.Range("A292").Copy()
$daily.Range("A297").PasteSpecial(-4163)
$daily.Range("A292").Copy()
$daily.Range("A298").PasteSpecial(-4163)

$daily.Range("B297").Value = "Workers Accommodation"
$daily.Range("C297").Value = "Food Supplies"
$daily.Range("D297").Value = "Unga and cooking oil for workers"
$daily.Range("E297").Value = 2700
$daily.Range("F297").Value = 25
$daily.Range("G297").Value = 2725
$daily.Range("H297").Value = "Supplier"
$daily.Range("I297").Value = "PAID"

$daily.Range("B298").Value = "Workers Accommodation"
$daily.Range("C298").Value = "Food"
$daily.Range("D298").Value = "Unga and cooking oil for workers"
$daily.Range("E298").Value = 2700
$daily.Range("F298").Value = 25
$daily.Range("G298").Value = 2725
$daily.Range("H298").Value = "Supplier"
$daily.Range("I298").Value = "PAID"

# Rows 291 & 296 just get the helper rate/price correction.
$daily.Range("D291").Value = "3 helpers @ 600 each - UNPAID"
$daily.Range("E291").Value = 1800
$daily.Range("D296").Value = "3 helpers @ 600 each - UNPAID"
$daily.Range("E296").Value = 1800

# ---------------------------------------------------------------------------
# Sheet: M-Pesa Fees
# ---------------------------------------------------------------------------
$mpesa = $wb.Worksheets.Item("M-Pesa Fees")

$mpesa.Range("C5").Value  = 37
$mpesa.Range("D5").Value  = 775
$mpesa.Range("B20").Value = "KES 9,963"

# ---------------------------------------------------------------------------
# Sheet: Unpaid Labor
# ---------------------------------------------------------------------------
$unpaid = $wb.Worksheets.Item("Unpaid Labor")

$unpaid.Range("B25").Value = "3 helpers @ 600 each - UNPAID"
$unpaid.Range("C25").Value = 1800
$unpaid.Range("B30").Value = "3 helpers @ 600 each - UNPAID"
$unpaid.Range("C30").Value = 1800
$unpaid.Range("C37").Value = "KES 46,100"
